$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'27.567.78"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "'1.846.66"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -1.27%  "
$ws.Range("D5").Value = "'333.90"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  -1.09%  "
$ws.Range("D7").Value = "'0.4655"
$ws.Range("E7").Value = "  -0.88%  "
$ws.Range("D8").Value = "'0.3862"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").Value = "'46.20"
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("D10").Value = "'0.07924"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").Value = "'0.9956"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").Value = "'21.50"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").Value = "'1.847.76"
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").Value = "'5.931"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").Value = "'7.123"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").Value = "'1.006"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "'89.02"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").Value = "'0.06661"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "'1.004"
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("D22").Value = "'27.570.23"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").Value = "'5.388"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("D24").Value = "'10.93"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").Value = "'2.308"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("D26").Value = "'158.07"
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("D27").Value = "'19.54"
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("D28").Value = "'2.105"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").Value = "'5.410"
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("D30").Value = "'119.94"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").Value = "'0.9775"
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").Value = "'0.09412"
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("D33").Value = "'3.589"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("D34").Value = "'5.291"
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("D35").Value = "'1.345"
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("D36").Value = "'0.06042"
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("D37").Value = "'0.02232"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("D38").Value = "'8.328"
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("D39").Value = "'1.183"
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("D40").Value = "'0.5896"
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").Value = "'0.1866"
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("D42").Value = "'10.34"
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").Value = "'1.244"
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("D44").Value = "'0.5589"
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("D45").Value = "'12.24"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("D47").Value = "'0.06690"
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D49").Value = "'1.053"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("D50").Value = "'1.005"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("E51").Value = "  -1.08%  "
